$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.543.62'
$ws.Range('E2').Value = '  -0.71%  '
$ws.Range('D3').Value = '1.596.43'
$ws.Range('E3').Value = '  -0.95%  '
$ws.Range('E4').Value = '  +0.67%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '208.30'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.09%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.499'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -3.81%  '
$ws.Range('E7').Value = '  +0.72%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '22.24'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -3.95%  '
$ws.Range('E9').Value = '  -1.43%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0588'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -3.00%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0870'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.40%  '
$ws.Range('D12').Value = '1.823.98'
$ws.Range('E12').Value = '  -0.61%  '
$ws.Range('D13').Value = '1.607.70'
$ws.Range('E13').Value = '  +0.09%  '
$ws.Range('E14').Value = '  -3.12%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.540'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -2.96%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '63.40'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -2.01%  '
$ws.Range('D17').Value = '27.536.23'
$ws.Range('E17').Value = '  -0.43%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '216.68'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -4.86%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.39'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -2.09%  '
$ws.Range('D20').Value = '0.0₃0689'
$ws.Range('E20').Value = '  -3.70%  '
$ws.Range('E21').Value = '  +0.87%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.19'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -1.76%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '9.74'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -2.96%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.01'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.40%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '154.55'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.45%  '
$ws.Range('E26').Value = '  +0.78%  '
$ws.Range('E27').Value = '  -2.02%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '15.01'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -2.37%  '
$ws.Range('E29').Value = '  -4.20%  '
$ws.Range('E30').Value = '  +0.03%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.0466'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -2.24%  '
$ws.Range('E32').Value = '  -2.22%  '
$ws.Range('D33').Value = '1.362.00'
$ws.Range('E33').Value = '  -1.39%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.95'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -3.48%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.54'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -1.01%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.961'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -2.52%  '
$ws.Range('E37').Value = '  -0.02%  '
$ws.Range('E38').Value = '  -1.84%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.536'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -2.91%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.811'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -3.96%  '
$ws.Range('E41').Value = '  +0.83%  '
$ws.Range('E42').Value = '  -4.69%  '
$ws.Range('B43').Value = 'Aave'
$ws.Range('C43').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '63.91'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -1.80%  '
$ws.Range('B44').Value = 'FraxShare'
$ws.Range('C44').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '5.31'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -1.75%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.75'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -3.23%  '
$ws.Range('D46').Value = '1.734.62'
$ws.Range('E46').Value = '  -0.59%  '
$ws.Range('E47').Value = '  -2.77%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '87.46'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.25%  '
$ws.Range('D49').Value = '0.0₇0995'
$ws.Range('E49').Value = '  +3.10%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0969'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -3.45%  '
$ws.Range('E51').Value = '  -0.83%  '
